# Update holdings weight/percent-change figures and the "as of" date in the
# confidential disclaimer note. The source sheet ships protected, so we briefly
# unprotect it to write the new values and then restore protection afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# "Model holdings provided as of 2021-03-24 ..." -> "...2021-03-25 ..."
$ws.Range("A58").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

# Refreshed Weight (col D) and Percent Change (col E) figures for rows 2-55
$ws.Range("D2").Value = 0.01663669925495738
$ws.Range("E2").Value = 0.01185578172368995
$ws.Range("D3").Value = 0.05148799300743596
$ws.Range("E3").Value = -0.01321965488310917
$ws.Range("D4").Value = 0.01485629250112356
$ws.Range("E4").Value = -0.006915598266609591
$ws.Range("D5").Value = 0.009797674368345449
$ws.Range("E5").Value = 0.01143946615824576
$ws.Range("D6").Value = 0.01574306008728758
$ws.Range("E6").Value = 0.005907355377854584
$ws.Range("D7").Value = 0.02078259790637001
$ws.Range("E7").Value = 0.006940253470126612
$ws.Range("D8").Value = 0.00440404995856443
$ws.Range("E8").Value = 0.009596521923545875
$ws.Range("D9").Value = 0.006765639064705493
$ws.Range("E9").Value = -0.00283990060347894
$ws.Range("D10").Value = 0.01403834080159272
$ws.Range("E10").Value = 0.01668056713928268
$ws.Range("D11").Value = 0.009051807560241799
$ws.Range("E11").Value = -0.00007370283018892643
$ws.Range("D12").Value = 0.01452238701144925
$ws.Range("E12").Value = 0.03996692392502732
$ws.Range("D13").Value = 0.002901241754687612
$ws.Range("E13").Value = 0.04024144869215296
$ws.Range("D14").Value = 0.006263479899714779
$ws.Range("E14").Value = -0.003514938488576425
$ws.Range("D15").Value = 0.01457035465322265
$ws.Range("E15").Value = 0.01281370335944776
$ws.Range("D16").Value = 0.01051952399837063
$ws.Range("E16").Value = 0.02340182648401834
$ws.Range("D17").Value = 0.02211231564212296
$ws.Range("E17").Value = 0.009300115402891995
$ws.Range("D18").Value = 0.008728242800017941
$ws.Range("E18").Value = 0.00786516853932584
$ws.Range("D19").Value = 0.01728276134526827
$ws.Range("E19").Value = 0.0003705762460626261
$ws.Range("D20").Value = 0.0120199304844396
$ws.Range("E20").Value = 0.01620691569073673
$ws.Range("D21").Value = 0.007402828142004057
$ws.Range("E21").Value = 0.003001000333444637
$ws.Range("D22").Value = 0.01367564805541602
$ws.Range("E22").Value = 0.01174702909438619
$ws.Range("D23").Value = 0.01984659510514513
$ws.Range("E23").Value = 0.01098202940642579
$ws.Range("D24").Value = 0.009879632988538233
$ws.Range("E24").Value = 0.02467443454420848
$ws.Range("D25").Value = 0.02087906690498521
$ws.Range("E25").Value = 0.02243083824873304
$ws.Range("D26").Value = 0.01120171192737455
$ws.Range("E26").Value = 0.0208867487961979
$ws.Range("D27").Value = 0.0202709986561068
$ws.Range("E27").Value = -0.008000724047425245
$ws.Range("D28").Value = 0.05608211224473809
$ws.Range("E28").Value = 0.004163544008660258
$ws.Range("D29").Value = 0.02027366723144884
$ws.Range("E29").Value = 0.02409937888198765
$ws.Range("D30").Value = 0.03058987914594685
$ws.Range("E30").Value = 0.003663962313530611
$ws.Range("D31").Value = 0.01569676030510297
$ws.Range("E31").Value = 0.0007862871520678549
$ws.Range("D32").Value = 0.013305716798624
$ws.Range("E32").Value = 0.02267314470801174
$ws.Range("D33").Value = 0.01993425780513155
$ws.Range("E33").Value = 0.0131325301204821
$ws.Range("D34").Value = 0.04067969579986349
$ws.Range("E34").Value = -0.00003443983606643553
$ws.Range("D35").Value = 0.01137613668316945
$ws.Range("E35").Value = 0.01090781140042218
$ws.Range("D36").Value = 0.009625818116317273
$ws.Range("E36").Value = 0.007956530176596033
$ws.Range("D37").Value = 0.01141066137665726
$ws.Range("E37").Value = 0.03069502302126725
$ws.Range("D38").Value = 0.007362265796804861
$ws.Range("E38").Value = 0.02419464455620512
$ws.Range("D39").Value = 0.0114471875016516
$ws.Range("E39").Value = 0.02059620596205947
$ws.Range("D40").Value = 0.01755555645960958
$ws.Range("E40").Value = 0.01732124874118823
$ws.Range("D41").Value = 0.0170425228501
$ws.Range("E41").Value = 0.01018966158423207
$ws.Range("D42").Value = 0.03330882384758051
$ws.Range("E42").Value = -0.003875619648490369
$ws.Range("D43").Value = 0.01125100585109294
$ws.Range("E43").Value = 0.01216536240341948
$ws.Range("D44").Value = 0.02151595576755779
$ws.Range("E44").Value = -0.0004806074878647149
$ws.Range("D45").Value = 0.01372791877492847
$ws.Range("E45").Value = 0.00968793054431738
$ws.Range("D46").Value = 0.008054460883341483
$ws.Range("E46").Value = 0.02195385590219523
$ws.Range("D47").Value = 0.01328330076575076
$ws.Range("E47").Value = 0.01365095149844064
$ws.Range("D48").Value = 0.009795406079304703
$ws.Range("E48").Value = 0.02602059580728211
$ws.Range("D49").Value = 0.01446534621351288
$ws.Range("E49").Value = 0.006572120373573132
$ws.Range("D50").Value = 0.008196462448730447
$ws.Range("E50").Value = 0.01418292514182951
$ws.Range("D51").Value = 0.01059184239014025
$ws.Range("E51").Value = 0.02979894686452855
$ws.Range("D52").Value = 0.008871245081160177
$ws.Range("E52").Value = -0.0001654464987386195
$ws.Range("D53").Value = 0.1449103125118014
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 0.04400480739044434
$ws.Range("E54").Value = 0.001546391752577536
$ws.Range("D55").Value = 0.9999999999999999
$ws.Range("E55").Value = 0.006976097593344166

# Restore sheet protection to its original (protected) state
$ws.Protect("ao6ra25", $true, $true, $true, $false, $true, $false, $false)
